# Remove the "district heating - bio gas-fired boiler" and
# "district heating - agricultural bio gas-fired boiler" rows from the
# HEATING sheet (commit: "removing agricultural bio gas fired boiler").
#
# In the original HEATING sheet:
#   row 5 = district heating - bio gas-fired boiler / T23
#   row 6 = district heating - agricultural bio gas-fired boiler / T24
#   row 7 = district heating - natural gas-fired boiler / T25
#
# Deleting rows 5 and 6 (bottom-up, so row indices stay valid) shifts the
# natural-gas-fired boiler row up to become the new row 5, matching the
# target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HEATING")

[void]$ws.Rows.Item(6).Delete()
[void]$ws.Rows.Item(5).Delete()

# Restore the view/selection state recorded in the target workbook: COOLING
# was visited with E20 selected, and HEATING ends up the active sheet with
# A13 selected.
$wsCooling = $wb.Worksheets.Item("COOLING")
[void]$wsCooling.Activate()
[void]$wsCooling.Range("E20").Select()

[void]$ws.Activate()
[void]$ws.Range("A13").Select()
